$d = $word.ActiveDocument

# --- Change 1: split "present in admin large parcels " into two runs ---
$p8 = $d.Paragraphs(8)
$paraText = $p8.Range.Text
$full = $p8.Range
$needle = "present in admin large parcels "
$splitAt = $full.Start + $paraText.IndexOf($needle)
$splitEnd = $splitAt + $needle.Length
$rngPresent = $d.Range($splitAt, $splitEnd)
$rngPresent.Text = "present in admin large parcels"
$rngPresent.Collapse(0)
$rngPresent.InsertAfter(". ")

# --- Change 2: insert two new sub-bullet paragraphs after the Chugach SP paragraph ---
$p8 = $d.Paragraphs(8)
$null = $p8.Range.InsertParagraphAfter()
$pAlt = $d.Paragraphs(9)
$xmlAlt = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:rPr><w:highlight w:val="yellow"/></w:rPr></w:pPr><w:r><w:t xml:space="preserve">Alternatively, </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>isn’t</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> there a service for this? – look at </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>axo</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> message from </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>sp</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> guy. </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$pAlt.Range.InsertXML($xmlAlt)

$pAltAfter = $d.Paragraphs(9)
$null = $pAltAfter.Range.InsertParagraphAfter()
$pConsider = $d.Paragraphs(10)
$xmlConsider = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr><w:rPr><w:highlight w:val="yellow"/></w:rPr></w:pPr><w:r><w:t xml:space="preserve">Consider updating to the service. </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$pConsider.Range.InsertXML($xmlConsider)

# --- Change 3: rewrite the "Arctic circle" paragraph (now paragraph 11) ---
$pArctic = $d.Paragraphs(11)
$xmlArctic = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:highlight w:val="yellow"/></w:rPr></w:pPr><w:r><w:rPr><w:strike/><w:highlight w:val="yellow"/></w:rPr><w:t xml:space="preserve">Arctic </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:strike/><w:highlight w:val="yellow"/></w:rPr><w:t xml:space="preserve">circle </w:t></w:r><w:r><w:rPr><w:strike/><w:highlight w:val="yellow"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:strike/></w:rPr><w:t>-</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> Not sure</w:t></w:r><w:r><w:t xml:space="preserve"> why new source is causing errors. May be simplest manually update for now. Alternatively load new version</w:t></w:r><w:r><w:t xml:space="preserve"> to</w:t></w:r><w:r><w:t xml:space="preserve"> common </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>sde</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> from stow</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$pArctic.Range.InsertXML($xmlArctic)

Write-Host "Done."
